# Actualización automática 2025-09-08 16:05:08
#
# Updates sales figures for client "CERAMIKASA S.A.S." (row 14 on the
# per-group / per-month sheets) and propagates the resulting totals into
# the monthly-total row and the monthly-compliance summary sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": new sales recorded for CERAMIKASA S.A.S.
# in the INODOROS (H), LAVABOS (I) and PIEDRA SINTERIZADA (L) columns.
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("H14").Value = 798.3
$wsGrupo.Range("I14").Value = 172.8
$wsGrupo.Range("L14").Value = 556.8099999999999

# Row 55 keeps a running "N de 53" count of clients with sales > 0 for
# each product column; it moves in lock-step with the H/I/L edits above.
$wsGrupo.Range("H55").Value = "1 de 53"
$wsGrupo.Range("I55").Value = "4 de 53"
$wsGrupo.Range("L55").Value = "3 de 53"

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL": same client's september ("septiembre") column
# and the grand-total row at the bottom of the sheet.
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F14").Value = 1858.55
$wsMensual.Range("F59").Value = 12182.43

# ---------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL": VENTA (D), POR CUMPLIR (E) and
# CUMPLIMIENTO (F) recompute for the affected product groups plus the
# TOTAL row.
# ---------------------------------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# INODOROS (row 6)
$wsCumpl.Range("D6").Value = 798.3
$wsCumpl.Range("E6").Value = 2109.28368146026
$wsCumpl.Range("F6").Value = 0.2745578760433385

# LAVABOS (row 7)
$wsCumpl.Range("D7").Value = 564.3
$wsCumpl.Range("E7").Value = 322.4110162875741
$wsCumpl.Range("F7").Value = 0.6363967399013218

# PIEDRA SINTERIZADA (row 11)
$wsCumpl.Range("D11").Value = 3043.66
$wsCumpl.Range("E11").Value = 14787.7543984654
$wsCumpl.Range("F11").Value = 0.1706908903570735

# TOTAL (row 15)
$wsCumpl.Range("D15").Value = 11802.42
$wsCumpl.Range("E15").Value = 110252.4155108344
$wsCumpl.Range("F15").Value = 0.09669768469723876
